$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Initial Estimate (C), Actual Time (D), and Week 1 remaining (F) values
$ws.Range("C3").Value = 1.5
$ws.Range("F3").Value = 1.5

$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 3
$ws.Range("F5").Value = 3

$ws.Range("C6").Value = 2.5

$ws.Range("F7").Value = 3

$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("F9").Value = 1

$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1

$ws.Range("C12").Value = 4
$ws.Range("F12").Value = 4

$ws.Range("C13").Value = 2.5
$ws.Range("F13").Value = 2.5

$ws.Range("F14").Value = 1.5

$ws.Range("F15").Value = 1.5

$ws.Range("D17").Value = 2

# Move active selection to F3 as in the saved workbook
$ws.Range("F3").Select()

$wb.Save()
